$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Boston -> Kyoto)
$ws.Range("C2").Value = 353.7
$ws.Range("D2").Value = 3.75

# Row 3 (London -> Paris)
$ws.Range("C3").Value = 455.86
$ws.Range("D3").Value = 5.98
$ws.Range("F3").Value = 2.28

# Row 8 (Berlin -> Munich)
$ws.Range("E8").Value = 538.21
$ws.Range("F8").Value = 5.13

# Row 9 (San Francisco -> Seattle)
$ws.Range("E9").Value = 593.28
$ws.Range("F9").Value = 6.82
